# Fix matched/non-matched row counts in the summary report layout.
#
# For each of the two worksheets ("summary" and "revsummary") there are two
# comparison blocks. In each block the "Non-matching Rows" count (row 7 / 21)
# and the "Matching Rows" count (row 9 / 23), as well as the
# "(Source1) - Matching Rows" count (row 12 / 26), had been computed with the
# wrong numbers swapped/duplicated. Correct them here.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Block 1 (rows 3-13)
    $ws.Cells.Item(7, 2).Value  = 39.0     # B7
    $ws.Cells.Item(7, 8).Value  = 0.0      # H7
    $ws.Cells.Item(9, 2).Value  = 960.0    # B9
    $ws.Cells.Item(9, 8).Value  = 960.0    # H9
    $ws.Cells.Item(12, 2).Value = 960.0    # B12
    $ws.Cells.Item(12, 8).Value = 960.0    # H12

    # Block 2 (rows 17-27)
    $ws.Cells.Item(21, 2).Value = 7.0      # B21
    $ws.Cells.Item(21, 8).Value = 32.0     # H21
    $ws.Cells.Item(23, 2).Value = 960.0    # B23
    $ws.Cells.Item(23, 8).Value = 960.0    # H23
    $ws.Cells.Item(26, 2).Value = 960.0    # B26
    $ws.Cells.Item(26, 8).Value = 960.0    # H26
}
